$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (A269) down through the new rows so the
# new date cells inherit the same cell style (s="2") as the existing ones.
$ws.Range("A269").Copy()
$ws.Range("A270:A301").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Flat list of (row, colA, colB, colC, colD) quintuples - nested arrays get
# flattened by this shell, so we index a flat array in strides of 5.
$data = @(
    270, 44344, 0, 2, 23.5654530458348,
    271, 44345, 0, 0, 0,
    272, 44346, 0, 0, 0,
    273, 44347, 1, 1, 11.7827265229174,
    274, 44348, 0, 1, 11.7827265229174,
    275, 44349, 0, 1, 11.7827265229174,
    276, 44350, 0, 1, 11.7827265229174,
    277, 44351, 0, 1, 11.7827265229174,
    278, 44352, 0, 1, 11.7827265229174,
    279, 44353, 0, 1, 11.7827265229174,
    280, 44354, 1, 1, 11.7827265229174,
    281, 44355, 0, 1, 11.7827265229174,
    282, 44356, 0, 1, 11.7827265229174,
    283, 44357, 0, 1, 11.7827265229174,
    284, 44358, 0, 1, 11.7827265229174,
    285, 44359, 0, 1, 11.7827265229174,
    286, 44360, 0, 1, 11.7827265229174,
    287, 44361, 0, 0, 0,
    288, 44362, 1, 1, 11.7827265229174,
    289, 44363, 0, 1, 11.7827265229174,
    290, 44364, 0, 1, 11.7827265229174,
    291, 44365, 1, 2, 23.5654530458348,
    292, 44366, 0, 2, 23.5654530458348,
    293, 44367, 0, 2, 23.5654530458348,
    294, 44368, 1, 3, 35.34817956875221,
    295, 44369, 0, 2, 23.5654530458348,
    296, 44370, 0, 2, 23.5654530458348,
    297, 44371, 4, 6, 70.69635913750442,
    298, 44372, 0, 5, 58.91363261458702,
    299, 44373, 1, 6, 70.69635913750442,
    300, 44374, 0, 6, 70.69635913750442,
    301, 44375, 0, 5, 58.91363261458702
)

for ($i = 0; $i -lt $data.Count; $i += 5) {
    $row = $data[$i]
    $ws.Cells.Item($row, 1).Value = $data[$i + 1]
    $ws.Cells.Item($row, 2).Value = $data[$i + 2]
    $ws.Cells.Item($row, 3).Value = $data[$i + 3]
    $ws.Cells.Item($row, 4).Value = $data[$i + 4]
}
